$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A7:F8")
$dataRange.NumberFormat = "@"

$ws.Range("A7").Value = "9691387881182728685"
$ws.Range("B7").Value = "28"
$ws.Range("C7").Value = "10"
$ws.Range("D7").Value = "2022"
$ws.Range("E7").Value = "Just a friendly reminder to enjoy your day"
$ws.Range("F7").Value = "1"

$ws.Range("A8").Value = "2880808233611366893"
$ws.Range("B8").Value = "10"
$ws.Range("C8").Value = "10"
$ws.Range("D8").Value = "2023"
$ws.Range("E8").Value = "Test Event (Again)"
$ws.Range("F8").Value = "0"

$dataRange.Style = "Normal"
